# "Add files via upload" -- appends two new example rows to the
# "code-list" sheet (rows 56 and 57), continuing the EXAMPLE# counter in
# column A, CATEGORY in column B and DESCRIPTION in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("code-list")

# Column A keeps incrementing from the row above (=A55+1, =A56+1, ...).
$ws.Range("A56").Formula = "=A55+1"
$ws.Range("A57").Formula = "=A56+1"

# Row 56: php/XML example - deleting an XML record.
$ws.Range("B56").Value = "php/XML"
$ws.Range("C56").Value = "Utility that deletes XML record"

# Row 57: JavaScript/JSON example - dynamic checkbox in a table.
$ws.Range("B57").Value = "JavaScript/JSON"
$ws.Range("C57").Value = "Add Checkbox to a dynamically created table"

# Leave the selection on the last cell typed into, like the author did.
$null = $ws.Range("D57").Select()
